$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[49.96168362564926, 50.036520579275276]"
$ws.Range("T2").Value = "[49.963202153099765, 50.015952692991185]"
$ws.Range("L3").Value = "[49.97892874539043, 50.11559541330282]"
$ws.Range("T3").Value = "[49.95639116955363, 50.036359056685775]"
